$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data row (row 78) with the same shape as the existing rows.
$ws.Cells.Item(78, 1).Value = "'2026/01/27"
$ws.Cells.Item(78, 2).Value = "逃离鸭科夫"
$ws.Cells.Item(78, 3).Value = 1160

# Copy the formatting from the previous row (77) so the new row matches
# the existing style (centered alignment, etc.).
$ws.Range("A77:C77").Copy()
$ws.Range("A78:C78").PasteSpecial(-4122)
